# Refresh the KHL referees stats workbook:
#   - "Главные" sheet (index 2) and "Линейные" sheet (index 3) get updated
#     per-referee stat totals for the rows whose underlying game counts
#     changed (one additional game recorded per affected referee), and
#   - every data row (2-26) on both sheets gets its "as_of_utc" (column AA)
#     refresh timestamp bumped to the new scrape time.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-13 03:04:09"
$aaCol = 27  # column AA

# --- Sheet "Главные" (index 2) ---------------------------------------------
$wsMain = $wb.Worksheets.Item(2)

$mainUpdates = @{
    5  = @{ 3=24;  4=399; 5=210; 6=189; 7=16.63; 8=8.75;  9=7.88;  10=100; 11=87; 12=2; 13=3; 22=18 }
    8  = @{ 3=22;  4=402; 5=201; 6=201; 7=18.27; 8=9.14;  9=9.14;  10=93;  11=93;             22=10 }
    11 = @{ 3=17;  4=455; 5=209; 6=246; 7=26.76; 8=12.29; 9=14.47; 10=92;  11=78; 12=3; 13=8;  22=12 }
    12 = @{ 3=15;  4=262; 5=106; 6=156; 7=17.47; 8=7.07;  9=10.4;  10=43;  11=53;        13=6 }
    14 = @{ 3=16;  4=192; 5=100; 6=92;  7=12;    8=6.25;  9=5.75;  10=50;  11=41;        13=2 }
    16 = @{ 3=24;  4=468; 5=229; 6=239; 7=19.5;  8=9.54;  9=9.96;  10=87;  11=87 }
    19 = @{ 3=19;  4=336; 5=160; 6=176; 7=17.68; 8=8.42;  9=9.26;  10=75;  11=73;             22=10 }
    20 = @{ 3=23;  4=392; 5=162; 6=230; 7=17.04; 8=7.04;  9=10;    10=76;  11=85;       13=6 }
    21 = @{ 3=20;  4=288; 5=126; 6=162; 7=14.4;  8=6.3;   9=8.1;   10=53;  11=66;       13=6 }
    22 = @{ 3=18;  4=334; 5=134; 6=200; 7=18.56; 8=7.44;  9=11.11; 10=67;  11=70 }
}

foreach ($r in $mainUpdates.Keys) {
    $rowVals = $mainUpdates[$r]
    foreach ($c in $rowVals.Keys) {
        $wsMain.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}

for ($r = 2; $r -le 26; $r++) {
    $wsMain.Cells.Item($r, $aaCol).Value = $newTimestamp
}

# --- Sheet "Линейные" (index 3) --------------------------------------------
$wsLin = $wb.Worksheets.Item(3)

$linUpdates = @{
    3  = @{ 3=23;  4=330; 5=162; 6=168; 7=14.35; 8=7.04; 9=7.3;   10=81; 11=69;              22=12 }
    11 = @{ 3=16;  4=226; 5=104; 6=122; 7=14.13; 8=6.5;  9=7.63;  10=52; 11=56;        13=2 }
    16 = @{ 3=23;  4=408; 5=193; 6=215; 7=17.74; 8=8.39; 9=9.35;  10=89; 11=95; 12=3; 13=5;  22=14 }
    22 = @{ 3=18;  4=287; 5=140; 6=147; 7=15.94; 8=7.78; 9=8.17;  10=70; 11=71 }
    24 = @{ 3=25;  4=456; 5=178; 6=278; 7=18.24; 8=7.12; 9=11.12; 10=79; 11=104; 12=2; 13=4; 22=10 }
    26 = @{ 3=21;  4=437; 5=183; 6=254; 7=20.81; 8=8.71; 9=12.1;  10=69; 11=72;        13=8 }
}

foreach ($r in $linUpdates.Keys) {
    $rowVals = $linUpdates[$r]
    foreach ($c in $rowVals.Keys) {
        $wsLin.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}

for ($r = 2; $r -le 26; $r++) {
    $wsLin.Cells.Item($r, $aaCol).Value = $newTimestamp
}

Write-Output "done"
